# Applies numeric-value corrections across the Tonberry_Profits leve-profit
# tracker sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), as produced by
# the scheduled Sheets runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 13889111
$ws.Range("I92").Value = 14706059
$ws.Range("K92").Value = 14706059
$ws.Range("M92").Value = -14704811

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1229.3334
$ws.Range("J97").Value = 1400
$ws.Range("L97").Value = 4200
$ws.Range("N97").Value = -5192

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 880.2143
$ws.Range("J129").Value = 899.2973
$ws.Range("L129").Value = 2697.8919
$ws.Range("N129").Value = -12697.8919

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6868
$ws.Range("I141").Value = 2721
$ws.Range("K141").Value = 8163
$ws.Range("M141").Value = -2983

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5752.0654
$ws.Range("I32").Value = 4029.158
$ws.Range("K32").Value = 4029.158
$ws.Range("M32").Value = -3742.158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1117
$ws.Range("I97").Value = 1069.8
$ws.Range("J97").Value = 1471
$ws.Range("K97").Value = 1069.8
$ws.Range("L97").Value = 1471
$ws.Range("M97").Value = -573.8
$ws.Range("N97").Value = -2463

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1305.6
$ws.Range("I102").Value = 1094.2
$ws.Range("K102").Value = 1094.2
$ws.Range("M102").Value = 527.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 922.1923
$ws.Range("I122").Value = 785.4091
$ws.Range("J122").Value = 1674.5
$ws.Range("K122").Value = 2356.2273
$ws.Range("L122").Value = 5023.5
$ws.Range("M122").Value = 93.77269999999999
$ws.Range("N122").Value = -9923.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1576.8
$ws.Range("I132").Value = 1279.2174
$ws.Range("K132").Value = 3837.6522
$ws.Range("M132").Value = -1307.6522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 54945
$ws.Range("J134").Value = 54945
$ws.Range("L134").Value = 54945
$ws.Range("N134").Value = -65085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 24874.445
$ws.Range("I11").Value = 1931.2
$ws.Range("J11").Value = 53553.5
$ws.Range("K11").Value = 1931.2
$ws.Range("L11").Value = 53553.5
$ws.Range("M11").Value = -1791.2
$ws.Range("N11").Value = -53833.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 338784.34
$ws.Range("I86").Value = 9068.666999999999
$ws.Range("J86").Value = 668500
$ws.Range("K86").Value = 9068.666999999999
$ws.Range("L86").Value = 668500
$ws.Range("M86").Value = -7945.666999999999
$ws.Range("N86").Value = -670746

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 338784.34
$ws.Range("I89").Value = 9068.666999999999
$ws.Range("J89").Value = 668500
$ws.Range("K89").Value = 45343.335
$ws.Range("L89").Value = 3342500
$ws.Range("M89").Value = -39727.335
$ws.Range("N89").Value = -3353732

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2072.3333
$ws.Range("J99").Value = 2072.3333
$ws.Range("L99").Value = 2072.3333
$ws.Range("N99").Value = -5068.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7100
$ws.Range("I134").Value = 8895.929
$ws.Range("J134").Value = 2909.5
$ws.Range("K134").Value = 26687.787
$ws.Range("L134").Value = 8728.5
$ws.Range("M134").Value = -24152.787
$ws.Range("N134").Value = -13798.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 200
$ws.Range("M7").Value = -87

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1557
$ws.Range("I31").Value = 1125.3077
$ws.Range("J31").Value = 2024.6666
$ws.Range("K31").Value = 1125.3077
$ws.Range("L31").Value = 2024.6666
$ws.Range("M31").Value = -830.3077000000001
$ws.Range("N31").Value = -2614.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1557
$ws.Range("I34").Value = 1125.3077
$ws.Range("J34").Value = 2024.6666
$ws.Range("K34").Value = 1125.3077
$ws.Range("L34").Value = 2024.6666
$ws.Range("M34").Value = -923.3077000000001
$ws.Range("N34").Value = -2428.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1652.5
$ws.Range("I58").Value = 870
$ws.Range("K58").Value = 870
$ws.Range("M58").Value = -667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4138.75
$ws.Range("I99").Value = 2185
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 2185
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -687
$ws.Range("N99").Value = -12996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 862.3077
$ws.Range("I105").Value = 787.7
$ws.Range("K105").Value = 787.7
$ws.Range("M105").Value = 959.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4138.75
$ws.Range("I126").Value = 2185
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 6555
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -4085
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1652.5
$ws.Range("I136").Value = 870
$ws.Range("K136").Value = 2610
$ws.Range("M136").Value = -60

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2775.3333
$ws.Range("J69").Value = 2985.4
$ws.Range("L69").Value = 8956.200000000001
$ws.Range("N69").Value = -10578.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2775.3333
$ws.Range("J72").Value = 2985.4
$ws.Range("L72").Value = 26868.6
$ws.Range("N72").Value = -34980.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 537
$ws.Range("J107").Value = 521.25
$ws.Range("L107").Value = 1563.75
$ws.Range("N107").Value = -5403.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 825.36365
$ws.Range("I121").Value = 864.5
$ws.Range("J121").Value = 816.6667
$ws.Range("K121").Value = 2593.5
$ws.Range("L121").Value = 2450.0001
$ws.Range("M121").Value = -1283.5
$ws.Range("N121").Value = -5070.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 28183.04
$ws.Range("J131").Value = 33188.137
$ws.Range("L131").Value = 99564.41100000001
$ws.Range("N131").Value = -109644.411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1128.8572
$ws.Range("I97").Value = 1128.8572
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1128.8572
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -632.8571999999999
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1048.6923
$ws.Range("I113").Value = 866.2
$ws.Range("K113").Value = 866.2
$ws.Range("M113").Value = 1303.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1648.9166
$ws.Range("I122").Value = 1597.375
$ws.Range("K122").Value = 4792.125
$ws.Range("M122").Value = -2342.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3973.2727
$ws.Range("I132").Value = 3761.6
$ws.Range("J132").Value = 4149.6665
$ws.Range("K132").Value = 11284.8
$ws.Range("L132").Value = 12448.9995
$ws.Range("M132").Value = -8754.799999999999
$ws.Range("N132").Value = -17508.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11933.467
$ws.Range("J40").Value = 8278.799999999999
$ws.Range("L40").Value = 8278.799999999999
$ws.Range("N40").Value = -8550.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2133.8096
$ws.Range("I61").Value = 1985
$ws.Range("J61").Value = 2375.625
$ws.Range("K61").Value = 1985
$ws.Range("L61").Value = 2375.625
$ws.Range("M61").Value = -1783
$ws.Range("N61").Value = -2779.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 19608814
$ws.Range("I93").Value = 864.3570999999999
$ws.Range("K93").Value = 864.3570999999999
$ws.Range("M93").Value = 383.6429000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2133.8096
$ws.Range("I113").Value = 1985
$ws.Range("J113").Value = 2375.625
$ws.Range("K113").Value = 1985
$ws.Range("L113").Value = 2375.625
$ws.Range("M113").Value = 185
$ws.Range("N113").Value = -6715.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5262.0835
$ws.Range("I122").Value = 4722.647
$ws.Range("J122").Value = 6572.143
$ws.Range("K122").Value = 14167.941
$ws.Range("L122").Value = 19716.429
$ws.Range("M122").Value = -11717.941
$ws.Range("N122").Value = -24616.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2303.361
$ws.Range("I132").Value = 1498.0769
$ws.Range("J132").Value = 2758.5217
$ws.Range("K132").Value = 4494.2307
$ws.Range("L132").Value = 8275.5651
$ws.Range("M132").Value = -1964.2307
$ws.Range("N132").Value = -13335.5651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74713.5
$ws.Range("J46").Value = 74713.5
$ws.Range("L46").Value = 74713.5
$ws.Range("N46").Value = -75175.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 4000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 4000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -6746

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1654.6538
$ws.Range("I132").Value = 1093.3636
$ws.Range("K132").Value = 3280.0908
$ws.Range("M132").Value = -750.0907999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 74713.5
$ws.Range("J134").Value = 74713.5
$ws.Range("L134").Value = 224140.5
$ws.Range("N134").Value = -229210.5

